$wb = $excel.ActiveWorkbook

# Sheet "展览" (Sheet1): rows 3-16
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 1067
$ws1.Range("F4").Value = 180
$ws1.Range("F5").Value = 2919
$ws1.Range("F6").Value = 96
$ws1.Range("F7").Value = 286
$ws1.Range("F8").Value = 30
$ws1.Range("F11").Value = 103
$ws1.Range("F12").Value = 151
$ws1.Range("F13").Value = 66
$ws1.Range("F14").Value = 2755
$ws1.Range("F15").Value = 1019
$ws1.Range("F16").Value = 3

# Sheet "全部类型" (Sheet4): rows 4-18
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value = 1067
$ws4.Range("F5").Value = 180
$ws4.Range("F6").Value = 2919
$ws4.Range("F7").Value = 96
$ws4.Range("F8").Value = 286
$ws4.Range("F9").Value = 30
$ws4.Range("F13").Value = 103
$ws4.Range("F14").Value = 151
$ws4.Range("F15").Value = 66
$ws4.Range("F16").Value = 2755
$ws4.Range("F17").Value = 1019
$ws4.Range("F18").Value = 3
